# Diagramas de flujo y ejecución actualizados
#
# 1) Bump the cached "datetimeFigureOut" field text from 2/6/22 to 3/6/22
#    everywhere it appears (slide master, every slide layout, notes master).
# 2) Re-flow / re-position the "ejecución" diagram shapes on slide 2 (the
#    unit-test row was removed, and the blocks below it slid up to close
#    the gap while the blocks above grew taller).
# 3) Remove the now-obsolete "unit_test" connector + label shapes.

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "2/6/22") {
                $shp.TextFrame.TextRange.Text = "3/6/22"
            }
        }
    }
}

# --- 1) Date placeholders -------------------------------------------------

Update-DatePlaceholders $p.SlideMaster.Shapes

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

Update-DatePlaceholders $p.NotesMaster.Shapes

# --- 2) Reposition / resize shapes on slide 2 ------------------------------

$s = $p.Slides.Item(2)

$shp = $s.Shapes.Item("CuadroTexto 62")
$shp.Top = 229.833740234375

$shp = $s.Shapes.Item("Conector recto de flecha 58")
$shp.Top = 107.95760345458984

$shp = $s.Shapes.Item("CuadroTexto 59")
$shp.Top = 78.2791748046875

$shp = $s.Shapes.Item("Conector recto de flecha 61")
$shp.Top = 229.4165802001953

$shp = $s.Shapes.Item("Conector recto de flecha 4")
$shp.Top = 157.721923828125

$shp = $s.Shapes.Item("Conector recto de flecha 29")
$shp.Top = 168.0627899169922

$shp = $s.Shapes.Item("CuadroTexto 30")
$shp.Top = 169.1616973876953

$shp = $s.Shapes.Item("CuadroTexto 33")
$shp.Top = 139.5438995361328

$shp = $s.Shapes.Item("Rectángulo 34")
$shp.Top = 192.73704528808594

$shp = $s.Shapes.Item("Conector recto de flecha 35")
$shp.Top = 216.5033416748047

$shp = $s.Shapes.Item("CuadroTexto 37")
$shp.Top = 197.93516540527344

$shp = $s.Shapes.Item("Conector recto de flecha 8")
$shp.Top = 302.8249206542969

$shp = $s.Shapes.Item("CuadroTexto 42")
$shp.Top = 303.7782287597656

$shp = $s.Shapes.Item("CuadroTexto 43")
$shp.Top = 121.67343139648438

$shp = $s.Shapes.Item("CuadroTexto 44")
$shp.Top = 193.23948669433594

$shp = $s.Shapes.Item("Conector recto de flecha 45")
$shp.Top = 346.5773010253906

$shp = $s.Shapes.Item("CuadroTexto 48")
$shp.Top = 328.4519348144531

$shp = $s.Shapes.Item("Rectángulo 54")
$shp.Top = 360.56036376953125

$shp = $s.Shapes.Item("CuadroTexto 55")
$shp.Top = 327.5906066894531

$shp = $s.Shapes.Item("CuadroTexto 60")
$shp.Top = 360.9006652832031

$shp = $s.Shapes.Item("Conector recto de flecha 14")
$shp.Top = 379.6166687011719

$shp = $s.Shapes.Item("CuadroTexto 86")
$shp.Top = 360.96588134765625

$shp = $s.Shapes.Item("Conector recto de flecha 16")
$shp.Top = 391.1754150390625

$shp = $s.Shapes.Item("CuadroTexto 88")
$shp.Top = 391.9447021484375

$shp = $s.Shapes.Item("Conector recto de flecha 18")
$shp.Top = 425.7328796386719

$shp = $s.Shapes.Item("Conector recto de flecha 20")
$shp.Top = 440.098388671875

$shp = $s.Shapes.Item("CuadroTexto 91")
$shp.Top = 407.2035827636719

$shp = $s.Shapes.Item("CuadroTexto 92")
$shp.Top = 441.0357971191406

$shp = $s.Shapes.Item("Conector recto de flecha 22")
$shp.Top = 488.8813171386719

$shp = $s.Shapes.Item("CuadroTexto 93")
$shp.Top = 489.8373718261719

$shp = $s.Shapes.Item("Rectángulo 2")
$shp.Left = 156.0437469482422
$shp.Top = 120.93484497070312
$shp.Width = 377.7890319824219

$shp = $s.Shapes.Item("Rectángulo 51")
$shp.Top = 327.2943115234375
$shp.Height = 182.11846923828125

# --- 3) Delete the obsolete "unit_test" arrow + label ----------------------

$s.Shapes.Item("Conector recto de flecha 52").Delete()
$s.Shapes.Item("CuadroTexto 53").Delete()
